$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 106.22581
$ws.Range("I33").Value = 102.85714
$ws.Range("K33").Value = 102.85714
$ws.Range("M33").Value = 126.14286
$ws.Range("H93").Value = 26633.334
$ws.Range("J93").Value = 26633.334
$ws.Range("L93").Value = 26633.334
$ws.Range("N93").Value = -31625.334
$ws.Range("H132").Value = 1276.2916
$ws.Range("I132").Value = 1323.4844
$ws.Range("J132").Value = 898.75
$ws.Range("K132").Value = 3970.4532
$ws.Range("L132").Value = 2696.25
$ws.Range("M132").Value = -1440.4532
$ws.Range("N132").Value = -7756.25
$ws.Range("H137").Value = 916.3461
$ws.Range("I137").Value = 883.7143
$ws.Range("K137").Value = 2651.1429
$ws.Range("M137").Value = -101.1428999999998

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1718
$ws.Range("I122").Value = 1225
$ws.Range("J122").Value = 1882.3334
$ws.Range("K122").Value = 3675
$ws.Range("L122").Value = 5647.0002
$ws.Range("M122").Value = -1225
$ws.Range("N122").Value = -10547.0002
$ws.Range("H132").Value = 1654.3334
$ws.Range("I132").Value = 1537.7142
$ws.Range("J132").Value = 2674.75
$ws.Range("K132").Value = 4613.142599999999
$ws.Range("L132").Value = 8024.25
$ws.Range("M132").Value = -2083.142599999999
$ws.Range("N132").Value = -13084.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 3680.5715
$ws.Range("I36").Value = 2000
$ws.Range("J36").Value = 4941
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 4941
$ws.Range("M36").Value = -1466
$ws.Range("N36").Value = -6009
$ws.Range("H134").Value = 16407.016
$ws.Range("I134").Value = 1458.1754
$ws.Range("K134").Value = 4374.5262
$ws.Range("M134").Value = -1839.5262

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1877364.5
$ws.Range("I31").Value = 3002627.8
$ws.Range("J31").Value = 1925.8334
$ws.Range("K31").Value = 3002627.8
$ws.Range("L31").Value = 1925.8334
$ws.Range("M31").Value = -3002332.8
$ws.Range("N31").Value = -2515.8334
$ws.Range("H34").Value = 1877364.5
$ws.Range("I34").Value = 3002627.8
$ws.Range("J34").Value = 1925.8334
$ws.Range("K34").Value = 3002627.8
$ws.Range("L34").Value = 1925.8334
$ws.Range("M34").Value = -3002425.8
$ws.Range("N34").Value = -2329.8334
$ws.Range("H58").Value = 3996.4324
$ws.Range("I58").Value = 1352.25
$ws.Range("J58").Value = 8878
$ws.Range("K58").Value = 1352.25
$ws.Range("L58").Value = 8878
$ws.Range("M58").Value = -1149.25
$ws.Range("N58").Value = -9284
$ws.Range("H122").Value = 1177.1666
$ws.Range("I122").Value = 624.4
$ws.Range("K122").Value = 1873.2
$ws.Range("M122").Value = 576.8000000000002
$ws.Range("H132").Value = 1490.2125
$ws.Range("I132").Value = 866.8125
$ws.Range("J132").Value = 2425.3125
$ws.Range("K132").Value = 2600.4375
$ws.Range("L132").Value = 7275.9375
$ws.Range("M132").Value = -70.4375
$ws.Range("N132").Value = -12335.9375
$ws.Range("H134").Value = 1198.541
$ws.Range("I134").Value = 1095.4889
$ws.Range("K134").Value = 3286.4667
$ws.Range("M134").Value = -751.4666999999999
$ws.Range("H136").Value = 3996.4324
$ws.Range("I136").Value = 1352.25
$ws.Range("J136").Value = 8878
$ws.Range("K136").Value = 4056.75
$ws.Range("L136").Value = 26634
$ws.Range("M136").Value = -1506.75
$ws.Range("N136").Value = -31734

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 484.4
$ws.Range("J122").Value = 484.8085
$ws.Range("L122").Value = 4363.2765
$ws.Range("N122").Value = -9263.2765

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3617.0833
$ws.Range("I80").Value = 3667.2222
$ws.Range("K80").Value = 3667.2222
$ws.Range("M80").Value = -2669.2222
$ws.Range("H83").Value = 3617.0833
$ws.Range("I83").Value = 3667.2222
$ws.Range("K83").Value = 18336.111
$ws.Range("M83").Value = -13344.111
$ws.Range("H97").Value = 1540
$ws.Range("I97").Value = 1540
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1540
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1044
$ws.Range("H102").Value = 1448.4445
$ws.Range("I102").Value = 1617.7142
$ws.Range("J102").Value = 856
$ws.Range("K102").Value = 1617.7142
$ws.Range("L102").Value = 856
$ws.Range("M102").Value = 4.285800000000108
$ws.Range("N102").Value = -4100
$ws.Range("H122").Value = 21638228
$ws.Range("I122").Value = 19956966
$ws.Range("J122").Value = 25000750
$ws.Range("K122").Value = 59870898
$ws.Range("L122").Value = 75002250
$ws.Range("M122").Value = -59868448
$ws.Range("N122").Value = -75007150
$ws.Range("H133").Value = 52450
$ws.Range("J133").Value = 52450
$ws.Range("L133").Value = 52450
$ws.Range("N133").Value = -62570

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4833580.5
$ws.Range("I7").Value = 3148.7693
$ws.Range("J7").Value = 11113142
$ws.Range("K7").Value = 3148.7693
$ws.Range("L7").Value = 11113142
$ws.Range("M7").Value = -3036.7693
$ws.Range("N7").Value = -11113366
$ws.Range("H40").Value = 2526953.2
$ws.Range("I40").Value = 5051404
$ws.Range("J40").Value = 2502.5
$ws.Range("K40").Value = 5051404
$ws.Range("L40").Value = 2502.5
$ws.Range("M40").Value = -5051268
$ws.Range("N40").Value = -2774.5
$ws.Range("H122").Value = 8827.4
$ws.Range("I122").Value = 11055.637
$ws.Range("J122").Value = 2699.75
$ws.Range("K122").Value = 33166.911
$ws.Range("L122").Value = 8099.25
$ws.Range("M122").Value = -30716.911
$ws.Range("N122").Value = -12999.25
$ws.Range("H126").Value = 4833580.5
$ws.Range("I126").Value = 3148.7693
$ws.Range("J126").Value = 11113142
$ws.Range("K126").Value = 9446.3079
$ws.Range("L126").Value = 33339426
$ws.Range("M126").Value = -6976.3079
$ws.Range("N126").Value = -33344366
$ws.Range("H132").Value = 1796.4127
$ws.Range("I132").Value = 1596.4286
$ws.Range("K132").Value = 4789.2858
$ws.Range("M132").Value = -2259.2858
$ws.Range("H136").Value = 1986.9149
$ws.Range("I136").Value = 1100.1282
$ws.Range("J136").Value = 6310
$ws.Range("K136").Value = 3300.3846
$ws.Range("L136").Value = 18930
$ws.Range("M136").Value = -750.3846000000003
$ws.Range("N136").Value = -24030

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1667
$ws.Range("I132").Value = 1703.2727
$ws.Range("J132").Value = 1622.6666
$ws.Range("K132").Value = 5109.8181
$ws.Range("L132").Value = 4867.9998
$ws.Range("M132").Value = -2579.8181
$ws.Range("N132").Value = -9927.9998
$ws.Range("H136").Value = 1763.5883
$ws.Range("I136").Value = 2013.9231
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 6041.7693
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = -3491.7693
$ws.Range("N136").Value = -7950

# --- GSM: remove N97 entirely (cell no longer present in edited row) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N97").ClearContents()
